$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking strings (e.g. "1.013", "27.197.10") that
# must stay as literal text, matching the source inlineStr cells. Pre-format
# the whole changed range as Text so Excel does not coerce them to numbers.
$ws.Range('D2:D51').NumberFormat = '@'

# Row 2
$ws.Range('D2').Value = '27.197.10'
$ws.Range('E2').Value = '  -3.53%  '

# Row 3
$ws.Range('D3').Value = '1.738.37'
$ws.Range('E3').Value = '  -3.29%  '

# Row 4
$ws.Range('D4').Value = '1.013'
$ws.Range('E4').Value = '  +0.83%  '

# Row 5
$ws.Range('D5').Value = '322.26'
$ws.Range('E5').Value = '  -4.88%  '

# Row 6
$ws.Range('D6').Value = '1.008'
$ws.Range('E6').Value = '  +0.72%  '

# Row 7
$ws.Range('D7').Value = '0.4193'
$ws.Range('E7').Value = '  -11.61%  '

# Row 8
$ws.Range('D8').Value = '0.3561'
$ws.Range('E8').Value = '  -5.08%  '

# Row 9
$ws.Range('D9').Value = '44.83'
$ws.Range('E9').Value = '  -1.55%  '

# Row 10
$ws.Range('D10').Value = '1.112'
$ws.Range('E10').Value = '  -2.93%  '

# Row 11
$ws.Range('D11').Value = '0.07314'
$ws.Range('E11').Value = '  -5.15%  '

# Row 12
$ws.Range('D12').Value = '1.011'
$ws.Range('E12').Value = '  +0.90%  '

# Row 13
$ws.Range('D13').Value = '21.31'
$ws.Range('E13').Value = '  -6.15%  '

# Row 14
$ws.Range('D14').Value = '6.053'
$ws.Range('E14').Value = '  -4.64%  '

# Row 15
$ws.Range('D15').Value = '7.154'
$ws.Range('E15').Value = '  -2.45%  '

# Row 16
$ws.Range('D16').Value = '1.750.59'
$ws.Range('E16').Value = '  -2.54%  '

# Row 17
$ws.Range('D17').Value = '0.00001049'
$ws.Range('E17').Value = '  -4.38%  '

# Row 18
$ws.Range('D18').Value = '84.60'
$ws.Range('E18').Value = '  +2.99%  '

# Row 19
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').Value = '0.05978'
$ws.Range('E19').Value = '  -11.29%  '

# Row 20
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = '1.009'
$ws.Range('E20').Value = '  +0.87%  '

# Row 21
$ws.Range('D21').Value = '16.68'
$ws.Range('E21').Value = '  -4.14%  '

# Row 22
$ws.Range('D22').Value = '6.020'
$ws.Range('E22').Value = '  -6.23%  '

# Row 23
$ws.Range('D23').Value = '27.278.45'
$ws.Range('E23').Value = '  -3.26%  '

# Row 24
$ws.Range('D24').Value = '11.22'
$ws.Range('E24').Value = '  -6.40%  '

# Row 25
$ws.Range('D25').Value = '2.414'
$ws.Range('E25').Value = '  +0.43%  '

# Row 26
$ws.Range('D26').Value = '19.78'
$ws.Range('E26').Value = '  -4.44%  '

# Row 27
$ws.Range('D27').Value = '149.27'
$ws.Range('E27').Value = '  -0.71%  '

# Row 28
$ws.Range('D28').Value = '2.306'
$ws.Range('E28').Value = '  -4.28%  '

# Row 29
$ws.Range('D29').Value = '1.949.36'
$ws.Range('E29').Value = '  -2.69%  '

# Row 30
$ws.Range('D30').Value = '1.277'
$ws.Range('E30').Value = '  +0.12%  '

# Row 31
$ws.Range('D31').Value = '126.33'
$ws.Range('E31').Value = '  -5.87%  '

# Row 32
$ws.Range('D32').Value = '3.740'
$ws.Range('E32').Value = '  -7.45%  '

# Row 33
$ws.Range('D33').Value = '0.08989'
$ws.Range('E33').Value = '  -7.08%  '

# Row 34
$ws.Range('D34').Value = '5.495'
$ws.Range('E34').Value = '  -7.57%  '

# Row 35
$ws.Range('D35').Value = '12.30'
$ws.Range('E35').Value = '  +0.88%  '

# Row 36
$ws.Range('D36').Value = '0.2131'
$ws.Range('E36').Value = '  -2.70%  '

# Row 37
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '0.02239'
$ws.Range('E37').Value = '  -5.80%  '

# Row 38
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = '0.06035'
$ws.Range('E38').Value = '  -4.42%  '

# Row 39
$ws.Range('D39').Value = '0.6385'
$ws.Range('E39').Value = '  -4.85%  '

# Row 40
$ws.Range('D40').Value = '4.947'
$ws.Range('E40').Value = '  -5.82%  '

# Row 41
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = '1.173'
$ws.Range('E41').Value = '  -3.93%  '

# Row 42
$ws.Range('B42').Value = 'Frax'
$ws.Range('C42').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D42').Value = '1.009'
$ws.Range('E42').Value = '  +0.88%  '

# Row 43
$ws.Range('D43').Value = '1.411'
$ws.Range('E43').Value = '  -4.91%  '

# Row 44
$ws.Range('D44').Value = '7.744'
$ws.Range('E44').Value = '  -4.74%  '

# Row 45
$ws.Range('D45').Value = '13.55'
$ws.Range('E45').Value = '  -4.23%  '

# Row 46
$ws.Range('D46').Value = '3.742'
$ws.Range('E46').Value = '  -3.51%  '

# Row 47
$ws.Range('D47').Value = '0.5824'
$ws.Range('E47').Value = '  -5.85%  '

# Row 48
$ws.Range('D48').Value = '123.36'
$ws.Range('E48').Value = '  -4.46%  '

# Row 49
$ws.Range('D49').Value = '1.920'
$ws.Range('E49').Value = '  -6.46%  '

# Row 50
$ws.Range('D50').Value = '0.06805'
$ws.Range('E50').Value = '  -4.14%  '

# Row 51
$ws.Range('D51').Value = '1.092'
$ws.Range('E51').Value = '  -7.01%  '
